$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing rows (14-21) that are no longer part of the table.
$ws.Range("A14:F21").EntireRow.Delete()

# Row 6: header row for the first time-entry block
$ws.Range("A6").Value = "Datum"
$ws.Range("B6").Value = "Bearbeiter"
$ws.Range("C6").Value = "Tätigkeitsbeschreibung"
$ws.Range("D6").Value = "Dauer"
$ws.Range("E6").Value = "Stunden-Satz"
$ws.Range("F6").Value = "Summe"

# Row 7: first time entry
$ws.Range("A7").Value = "19.09.2024"
$ws.Range("B7").Value = "Boytinck,Barbara (BEB)"
$ws.Range("C7").Value = "40551/2024 CHA-NI SE KG VAT reimbursement CHB Russia:Prüfung MSPA und LATA Russland; Telco mit C. Häußermann"
$ws.Range("D7").Value = "0,50"
$ws.Range("E7").Value = "400,00"
$ws.Range("F7").Value = "200,00"

# Row 8: subtotal for block 1
$ws.Range("A8").Value = "Gesamt"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "0,50"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "200,00"

# Row 9: header row for the second time-entry block
$ws.Range("A9").Value = "Datum"
$ws.Range("B9").Value = "Bearbeiter"
$ws.Range("C9").Value = "Tätigkeitsbeschreibung"
$ws.Range("D9").Value = "Dauer"
$ws.Range("E9").Value = "Stunden-Satz"
$ws.Range("F9").Value = "Summe"

# Row 10: second block, entry 1
$ws.Range("A10").Value = "25.09.2024"
$ws.Range("B10").Value = "Dr. Schlaffge,Andrea (ACS)"
$ws.Range("C10").Value = "30116/2024 Project Jura: Durchsicht Mail Frau Kues und PrüfungBehandlung Ambeo"
$ws.Range("D10").Value = "0,25"
$ws.Range("E10").Value = "500,00"
$ws.Range("F10").Value = "125,00"

# Row 11: second block, entry 2
$ws.Range("A11").Value = "26.09.2024"
$ws.Range("B11").Value = "Dr. Schlaffge,Andrea (ACS)"
$ws.Range("C11").Value = "30116/2024 Project Jura - AMBEO Soundbar Soundfiles: Durchsichtdes LATA hinsichtl. Behandlung der AMBEO sound files"
$ws.Range("D11").Value = "0,50"
$ws.Range("E11").Value = "500,00"
$ws.Range("F11").Value = "250,00"

# Row 12: second block, entry 3
$ws.Range("A12").Value = "26.09.2024"
$ws.Range("B12").Value = "Boytinck,Barbara (BEB)"
$ws.Range("C12").Value = "30116/2024 Project Jura - AMBEO Soundbar Soundfiles: PrüfungLATA SE KG; Alignment mit A. Schlaffge; Telco mit C. Häußermannund Frau Kues"
$ws.Range("D12").Value = "1,25"
$ws.Range("E12").Value = "400,00"
$ws.Range("F12").Value = "500,00"

# Row 13: subtotal for block 2
$ws.Range("A13").Value = "Gesamt"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = "2,00"
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = "875,00"

Write-Output "edit complete"
